# Update the "想去人数" (interested count) column F for several rows
# on both the "展览" sheet and the aggregated "全部类型" sheet.
# Changes:
#   F4:  60    -> 61
#   F5:  369   -> 371
#   F6:  11100 -> 11111
#   F7:  561   -> 568
#   F12: 149   -> 151
#   F19: 1182  -> 1185
#   F20: 57    -> 58

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F4"  = 61
    "F5"  = 371
    "F6"  = 11111
    "F7"  = 568
    "F12" = 151
    "F19" = 1185
    "F20" = 58
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
